# Update cryptocurrency "Price" column (column D) values to the latest
# snapshot pulled by the symbol-list GitHub Action.
# Cells store prices as plain text (not numbers), so each target cell is
# (re)formatted as Text ("@") before assignment to preserve exact digits
# (including trailing zeros and fixed decimal notation for very small values).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "245.61"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "23.15"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.407"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06048"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.8047"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9342"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1427"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07462"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.03368"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03072"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.010"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09364"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001604"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.04832"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0005941"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.005064"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.004162"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0009825"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.00008702"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.647"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.435"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.188"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03978"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006411"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1076"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002901"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.006299"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005261"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0005801"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.9002"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.002240"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002100"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.01010"
